# The two records (rows 7 and 8) in this sighting export were swapped:
# row 7 should hold the data that used to be in row 8, and vice versa,
# for every column where the two rows actually differ
# (A, B, E, F, G, H, I, Q, R, AC). The other columns already contain
# identical values in both rows, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "I", "Q", "R", "AC")

foreach ($col in $cols) {
    $cell7 = $ws.Range($col + "7")
    $cell8 = $ws.Range($col + "8")

    $v7 = $cell7.Value2
    $v8 = $cell8.Value2

    $cell7.Value2 = $v8
    $cell8.Value2 = $v7
}
